$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59; this shifts the existing rows 59-118
# down to 60-119 (and the dimension grows from R118 to R119).
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly data point.
$ws.Cells.Item(59, 1).Value = 6
$ws.Cells.Item(59, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(59, 3).Value = "Metropolitana"
$ws.Cells.Item(59, 4).Value = 45195
$ws.Cells.Item(59, 5).Value = 13
$ws.Cells.Item(59, 6).Value = 100112035
$ws.Cells.Item(59, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 220
$ws.Cells.Item(59, 11).Value = 17000
$ws.Cells.Item(59, 12).Value = 20000
$ws.Cells.Item(59, 13).Value = 18364
$ws.Cells.Item(59, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(59, 16).Value = 1224
$ws.Cells.Item(59, 17).Value = 15
$ws.Cells.Item(59, 18).Value = "Hortaliza"
